$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14, pushing the existing data (old rows 14-35) down to 15-36.
$ws.Rows.Item(14).Insert()

# Copy the date-format style from the row that got shifted (now row 15, column D)
# onto the new row's D cell so it keeps the same date number format.
$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the newly inserted row 14 with the new weekly record.
$ws.Cells.Item(14, 1).Value = 5
$ws.Cells.Item(14, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(14, 3).Value = "Maule"
$ws.Cells.Item(14, 4).Value = 44482
$ws.Cells.Item(14, 5).Value = 7
$ws.Cells.Item(14, 6).Value = 100112026
$ws.Cells.Item(14, 7).Value = "Haba"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 500
$ws.Cells.Item(14, 11).Value = 8000
$ws.Cells.Item(14, 12).Value = 8000
$ws.Cells.Item(14, 13).Value = 8000
$ws.Cells.Item(14, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(14, 16).Value = 320
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"
